$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- Step 1: introduce new shared strings in the exact order needed ---
$ws.Range("F151").Value = "Tests"
$ws.Range("G151").Value = "Tests de la strat" + [char]0x00E9 + "gie"
$ws.Range("G148").Value = "Manuel d'installation"
$ws.Range("G152").Value = "Sources, glossaire"
$ws.Range("G153").Value = "V" + [char]0x00E9 + "rification de l'exhaustivit" + [char]0x00E9
$ws.Range("H153").Value = "V" + [char]0x00E9 + "rification de l'exhaustivit" + [char]0x00E9 + " de la documentation."

# --- Step 2: fill in the rest of rows 148-154 ---

# Row 148
$ws.Range("A148").Value = 44350
$ws.Range("B148").Value = 5
$ws.Range("C148").Value = 0.33333333333333331
$ws.Range("D148").Value = 0.37847222222222227
$ws.Range("F148").Value = "Documentation"

# Row 149
$ws.Range("A149").Value = 44350
$ws.Range("B149").Value = 5
$ws.Range("C149").Value = 0.37847222222222227
$ws.Range("D149").Value = 0.39930555555555558
$ws.Range("F149").Value = "R" + [char]0x00E9 + "alisation"
$ws.Range("G149").Value = "Dossier de r" + [char]0x00E9 + "alisation"
$ws.Range("H149").Value = "Documentation"

# Row 150
$ws.Range("A150").Value = 44350
$ws.Range("B150").Value = 5
$ws.Range("C150").Value = 0.40972222222222227
$ws.Range("D150").Value = 0.4236111111111111
$ws.Range("F150").Value = "R" + [char]0x00E9 + "alisation"
$ws.Range("G150").Value = "Dossier de r" + [char]0x00E9 + "alisation"
$ws.Range("H150").Value = "Documentation"

# Row 151
$ws.Range("A151").Value = 44350
$ws.Range("B151").Value = 5
$ws.Range("C151").Value = 0.4236111111111111
$ws.Range("D151").Value = 0.4861111111111111

# Row 152
$ws.Range("A152").Value = 44350
$ws.Range("B152").Value = 5
$ws.Range("C152").Value = 0.4861111111111111
$ws.Range("D152").Value = 0.51041666666666663
$ws.Range("F152").Value = "Documentation"

# Row 153
$ws.Range("A153").Value = 44350
$ws.Range("B153").Value = 5
$ws.Range("C153").Value = 0.5625
$ws.Range("D153").Value = 0.62847222222222221
$ws.Range("F153").Value = "Documentation"

# Row 154
$ws.Range("A154").Value = 44350
$ws.Range("B154").Value = 5
$ws.Range("C154").Value = 0.63888888888888895
$ws.Range("F154").Value = "Documentation"
$ws.Range("G154").Value = $ws.Range("G153").Value()
$ws.Range("H154").Value = $ws.Range("H153").Value()

# row heights for the two wrapped rows
$ws.Rows.Item(153).RowHeight = 30
$ws.Rows.Item(154).RowHeight = 30

Write-Host "rows filled"

# --- Step 3: add new blank rows 155-161, copying the blank-row layout/style from row 151 ---
$ws.Range("A151:E151").Copy()
$ws.Range("A155:E161").PasteSpecial(-4122)
$ws.Range("G151:L151").Copy()
$ws.Range("G155:L161").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("E155:E161").Formula = "=D155-C155"

Write-Host "new rows added"

# --- Step 4: update the sheet view (scrolled position + active selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 139
$ws.Range("D154").Select()

Write-Host "view updated"
